# Generate Report for Handback
# Update the "Correspond Handback Datetime" (E) and "Correspond Handback
# DateTime" (H) values for the second data row (row 3) on the zh-cn and
# de-de sheets, reflecting a fresh handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 18:47:32"
$wsZhCn.Range("H3").Value = "2016-03-21 18:47:53"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 18:47:36"
$wsDeDe.Range("H3").Value = "2016-03-21 18:47:59"
